$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header labels in row 1
$ws.Range("B1").Value = "AVERAGE_MERGE_EFFORT_IN_100_COMMITS_WO_FT"
$ws.Range("C1").Value = "AVERAGE_MERGE_EFFORT_IN_100_COMMITS_WITH_FT"

# Update the active selection to B4
$ws.Range("B4").Select()
